# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions scheduled refresh).
#
# Column D ("Price") holds price strings as literal text (not numbers) so
# that values like "1.00" or "28.30" keep their trailing zeros and
# thousands-dot-separated values like "69.392.10" are not mangled by
# numeric parsing. We explicitly format those cells as Text ("@") before
# assigning a numeric-looking string so Excel stores them as text, matching
# the original workbook's inline-string cells.
# Column E ("Volume(1h)") values already contain spaces and a "%" sign, so
# they are safely preserved as text without any extra formatting step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.392.10"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "3.345.07"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "192.09"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.55"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.72"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.425"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").Value = "3.928.91"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.30"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "69.439.33"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "3.353.38"
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.83"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.75"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "427.22"
$ws.Range("E20").Value = "  +7.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.72"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.39"
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.63"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.03"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.28"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.02"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.38"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.03"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.810"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.58"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.746.61"
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.46"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.51"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.19"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.41"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0687"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "343.49"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.52"
$ws.Range("E49").Value = "  +4.32%  "
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.30"
$ws.Range("E51").Value = "  -0.31%  "
